$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People
$ws.Range("H32").Value = 95888.5
$ws.Range("I32").Value = 5952
$ws.Range("J32").Value = 221799.6
$ws.Range("K32").Value = 5952
$ws.Range("L32").Value = 221799.6
$ws.Range("M32").Value = -5626
$ws.Range("N32").Value = -222451.6

# Row 41: The Write Stuff
$ws.Range("H41").Value = 761.7619
$ws.Range("I41").Value = 530.8182
$ws.Range("K41").Value = 530.8182
$ws.Range("M41").Value = -90.81820000000005

# Row 64: Forged from the Void
$ws.Range("H64").Value = 3934.9807
$ws.Range("J64").Value = 4642.6665
$ws.Range("L64").Value = 4642.6665
$ws.Range("N64").Value = -5138.6665

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 3934.9807
$ws.Range("J67").Value = 4642.6665
$ws.Range("L67").Value = 4642.6665
$ws.Range("N67").Value = -6358.6665

# Row 103: Let Loose the Juice
$ws.Range("H103").Value = 181.8125
$ws.Range("I103").Value = 179.1
$ws.Range("J103").Value = 186.33333
$ws.Range("K103").Value = 537.3
$ws.Range("L103").Value = 558.99999
$ws.Range("M103").Value = 48.70000000000005
$ws.Range("N103").Value = -1730.99999

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 636.5
$ws.Range("I107").Value = 484.8
$ws.Range("J107").Value = 889.3333
$ws.Range("K107").Value = 484.8
$ws.Range("L107").Value = 889.3333
$ws.Range("M107").Value = 1435.2
$ws.Range("N107").Value = -4729.3333

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 1014.7742
$ws.Range("J112").Value = 998.9231
$ws.Range("L112").Value = 2996.7693
$ws.Range("N112").Value = -5212.7693

# Row 125: Body over Mind
$ws.Range("H125").Value = 2260.125
$ws.Range("I125").Value = 1095
$ws.Range("J125").Value = 2648.5
$ws.Range("K125").Value = 9855
$ws.Range("L125").Value = 23836.5
$ws.Range("M125").Value = -7395
$ws.Range("N125").Value = -28756.5

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2116.25
$ws.Range("I132").Value = 1898.1428
$ws.Range("K132").Value = 5694.428400000001
$ws.Range("M132").Value = -3164.428400000001

# Row 135: For Tired Minds
$ws.Range("H135").Value = 1315.8
$ws.Range("I135").Value = 804.5714
$ws.Range("J135").Value = 3999.75
$ws.Range("K135").Value = 7241.1426
$ws.Range("L135").Value = 35997.75
$ws.Range("M135").Value = -4706.1426
$ws.Range("N135").Value = -41067.75

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 470540.16
$ws.Range("I137").Value = 2625.8
$ws.Range("K137").Value = 7877.400000000001
$ws.Range("M137").Value = -5327.400000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7545.579
$ws.Range("I32").Value = 4098.548
$ws.Range("K32").Value = 4098.548
$ws.Range("M32").Value = -3811.548

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 20837354
$ws.Range("I45").Value = 5680.75
$ws.Range("J45").Value = 62500700
$ws.Range("K45").Value = 5680.75
$ws.Range("L45").Value = 62500700
$ws.Range("M45").Value = -5303.75
$ws.Range("N45").Value = -62501454

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1107.8431
$ws.Range("I132").Value = 927.0833
$ws.Range("K132").Value = 2781.2499
$ws.Range("M132").Value = -251.2498999999998

# Row 134: Brace for More Vambraces
$ws.Range("H134").Value = 108054.164
$ws.Range("J134").Value = 108054.164
$ws.Range("L134").Value = 108054.164
$ws.Range("N134").Value = -118194.164

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 126193.875
$ws.Range("I22").Value = 144199.86
$ws.Range("J22").Value = 152
$ws.Range("K22").Value = 144199.86
$ws.Range("L22").Value = 152
$ws.Range("M22").Value = -144026.86
$ws.Range("N22").Value = -498

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2455.842
$ws.Range("I31").Value = 2084.4075
$ws.Range("K31").Value = 2084.4075
$ws.Range("M31").Value = -1789.4075

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2455.842
$ws.Range("I34").Value = 2084.4075
$ws.Range("K34").Value = 2084.4075
$ws.Range("M34").Value = -1882.4075

# Row 95: Standing on Ceremony
$ws.Range("H95").Value = 10049
$ws.Range("I95").Value = 15500
$ws.Range("J95").Value = 9140.5
$ws.Range("K95").Value = 15500
$ws.Range("L95").Value = 9140.5
$ws.Range("M95").Value = -12754
$ws.Range("N95").Value = -14632.5

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1773.0769
$ws.Range("I132").Value = 1642.875
$ws.Range("J132").Value = 2368.2856
$ws.Range("K132").Value = 4928.625
$ws.Range("L132").Value = 7104.8568
$ws.Range("M132").Value = -2398.625
$ws.Range("N132").Value = -12164.8568

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 29327.553
$ws.Range("I134").Value = 2763.606
$ws.Range("K134").Value = 8290.818000000001
$ws.Range("M134").Value = -5755.818000000001

# Row 138: Bow Out
$ws.Range("H138").Value = 51490.9
$ws.Range("J138").Value = 49911.11
$ws.Range("L138").Value = 49911.11
$ws.Range("N138").Value = -60191.11

$ws = $wb.Worksheets.Item("CUL")
# Row 13: Fishy Revelations
$ws.Range("H13").Value = 500.25
$ws.Range("I13").Value = 500.5
$ws.Range("K13").Value = 1501.5
$ws.Range("M13").Value = -1333.5

# Row 15: Pretty Enough to Eat
$ws.Range("H15").Value = 41.363636
$ws.Range("I15").Value = 58
$ws.Range("J15").Value = 21.4
$ws.Range("K15").Value = 174
$ws.Range("L15").Value = 64.19999999999999
$ws.Range("M15").Value = -34
$ws.Range("N15").Value = -344.2

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 1300
$ws.Range("I55").Value = 1300
$ws.Range("K55").Value = 3900
$ws.Range("M55").Value = -3723

# Row 139: Najoothie
$ws.Range("H139").Value = 9143.120000000001
$ws.Range("I139").Value = 1920.7273
$ws.Range("K139").Value = 5762.1819
$ws.Range("M139").Value = -622.1818999999996

$ws = $wb.Worksheets.Item("GSM")
# Row 21: Forever 21K
$ws.Range("H21").Value = 306249.9
$ws.Range("I21").Value = 3400000
$ws.Range("K21").Value = 3400000
$ws.Range("M21").Value = -3399827

# Row 30: Dog Tags Are for Dogs
$ws.Range("H30").Value = 306249.9
$ws.Range("I30").Value = 3400000
$ws.Range("K30").Value = 3400000
$ws.Range("M30").Value = -3399895

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 1399.0476
$ws.Range("I102").Value = 1399.0476
$ws.Range("K102").Value = 1399.0476
$ws.Range("M102").Value = 222.9523999999999

# Row 132: On Board for Lar
$ws.Range("H132").Value = 7709.25
$ws.Range("I132").Value = 6503
$ws.Range("K132").Value = 19509
$ws.Range("M132").Value = -16979

# Row 134: Guaranteed Gem
$ws.Range("H134").Value = 53704.168
$ws.Range("J134").Value = 53704.168
$ws.Range("L134").Value = 161112.504
$ws.Range("N134").Value = -166182.504

$ws = $wb.Worksheets.Item("LTW")
# Row 9: From the Sands to the Stage
$ws.Range("H9").Value = 7644.0713
$ws.Range("J9").Value = 25624
$ws.Range("L9").Value = 25624
$ws.Range("N9").Value = -26072

# Row 13: Throwing Down the Gauntlet
$ws.Range("H13").Value = 7495
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 7495
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 7495
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -7775

# Row 17: Only the Best
$ws.Range("H17").Value = 10549.5
$ws.Range("J17").Value = 13749
$ws.Range("L17").Value = 13749
$ws.Range("N17").Value = -14089

# Row 19: Targe Up
$ws.Range("H19").Value = 7899.8
$ws.Range("J19").Value = 7899.8
$ws.Range("L19").Value = 7899.8
$ws.Range("N19").Value = -8239.799999999999

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 8680.691999999999
$ws.Range("I22").Value = 1114.1428
$ws.Range("J22").Value = 17508.334
$ws.Range("K22").Value = 1114.1428
$ws.Range("L22").Value = 17508.334
$ws.Range("M22").Value = -819.1428000000001
$ws.Range("N22").Value = -18098.334

# Row 27: Fire and Hide
$ws.Range("H27").Value = 8680.691999999999
$ws.Range("I27").Value = 1114.1428
$ws.Range("J27").Value = 17508.334
$ws.Range("K27").Value = 1114.1428
$ws.Range("L27").Value = 17508.334
$ws.Range("M27").Value = -1007.1428
$ws.Range("N27").Value = -17722.334

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 14431
$ws.Range("I46").Value = 18574.834
$ws.Range("J46").Value = 1999.5
$ws.Range("K46").Value = 18574.834
$ws.Range("L46").Value = 1999.5
$ws.Range("M46").Value = -18386.834
$ws.Range("N46").Value = -2375.5

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 2662.4666
$ws.Range("I61").Value = 2424.1428
$ws.Range("K61").Value = 2424.1428
$ws.Range("M61").Value = -2222.1428

# Row 113: Peace in Rest
$ws.Range("H113").Value = 2662.4666
$ws.Range("I113").Value = 2424.1428
$ws.Range("K113").Value = 2424.1428
$ws.Range("M113").Value = -254.1428000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 1547.909
$ws.Range("I81").Value = 881.125
$ws.Range("J81").Value = 3326
$ws.Range("K81").Value = 1762.25
$ws.Range("L81").Value = 6652
$ws.Range("M81").Value = -701.25
$ws.Range("N81").Value = -8774

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 1547.909
$ws.Range("I84").Value = 881.125
$ws.Range("J84").Value = 3326
$ws.Range("K84").Value = 8811.25
$ws.Range("L84").Value = 33260
$ws.Range("M84").Value = -3507.25
$ws.Range("N84").Value = -43868

# Row 93: What Guides Want
$ws.Range("H93").Value = 82000.336
$ws.Range("J93").Value = 75000
$ws.Range("L93").Value = 75000
$ws.Range("N93").Value = -79992

# Row 113: A Tender Table
$ws.Range("H113").Value = 1100.8572
$ws.Range("I113").Value = 1336.8334
$ws.Range("J113").Value = 786.2222
$ws.Range("K113").Value = 4010.5002
$ws.Range("L113").Value = 2358.6666
$ws.Range("M113").Value = -1840.5002
$ws.Range("N113").Value = -6698.6666

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 34837.97
$ws.Range("I126").Value = 40458.5
$ws.Range("K126").Value = 121375.5
$ws.Range("M126").Value = -118905.5

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 967663.9
$ws.Range("I132").Value = 1431.2858
$ws.Range("K132").Value = 4293.857400000001
$ws.Range("M132").Value = -1763.857400000001

# Row 133: Begin with the Basics
$ws.Range("H133").Value = 66637.60000000001
$ws.Range("J133").Value = 64422
$ws.Range("L133").Value = 64422
$ws.Range("N133").Value = -74542

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1231.7778
$ws.Range("I136").Value = 747.7778
$ws.Range("J136").Value = 3167.7778
$ws.Range("K136").Value = 2243.3334
$ws.Range("L136").Value = 9503.3334
$ws.Range("M136").Value = 306.6666
$ws.Range("N136").Value = -14603.3334

# Row 137: Traditional Trousers
$ws.Range("H137").Value = 149998.33
$ws.Range("J137").Value = 149998.33
$ws.Range("L137").Value = 149998.33
$ws.Range("N137").Value = -160198.33
